$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.393.68'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").Value = '1.859.93'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '324.52'
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").Value = '0.4549'
$ws.Range("E7").Value = '  -2.10%  '

$ws.Range("D8").Value = '0.3823'
$ws.Range("E8").Value = '  -1.88%  '

$ws.Range("D9").Value = '0.07801'
$ws.Range("E9").Value = '  -1.06%  '

$ws.Range("D10").Value = '0.9830'
$ws.Range("E10").Value = '  +1.07%  '

$ws.Range("D11").Value = '21.47'
$ws.Range("E11").Value = '  -3.35%  '

$ws.Range("D12").Value = '1.861.11'
$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("D13").Value = '5.625'
$ws.Range("E13").Value = '  -1.59%  '

$ws.Range("D14").Value = '6.876'
$ws.Range("E14").Value = '  -0.83%  '

$ws.Range("D15").Value = '0.06910'
$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  +0.35%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '86.46'
$ws.Range("E17").Value = '  -2.70%  '

$ws.Range("D18").Value = '0.000009931'
$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("D19").Value = '16.64'
$ws.Range("E19").Value = '  -1.25%  '

$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").Value = '28.400.96'
$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").Value = '5.246'
$ws.Range("E22").Value = '  -1.50%  '

$ws.Range("D23").Value = '10.85'
$ws.Range("E23").Value = '  -1.79%  '

$ws.Range("D24").Value = '2.096'
$ws.Range("E24").Value = '  -1.25%  '

$ws.Range("D25").Value = '2.069.22'
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").Value = '153.54'
$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("D27").Value = '19.06'
$ws.Range("E27").Value = '  -1.06%  '

$ws.Range("D28").Value = '5.616'
$ws.Range("E28").Value = '  -2.70%  '

$ws.Range("D29").Value = '117.21'
$ws.Range("E29").Value = '  -1.55%  '

$ws.Range("D30").Value = '1.903'
$ws.Range("E30").Value = '  -4.31%  '

$ws.Range("D31").Value = '0.09279'
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("D32").Value = '0.9031'
$ws.Range("E32").Value = '  -3.60%  '

$ws.Range("D33").Value = '5.251'
$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("D34").Value = '1.310'
$ws.Range("E34").Value = '  -1.93%  '

$ws.Range("D35").Value = '3.300'
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").Value = '0.05668'
$ws.Range("E36").Value = '  -2.91%  '

$ws.Range("D37").Value = '1.151'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").Value = '0.02046'
$ws.Range("E38").Value = '  -3.32%  '

$ws.Range("D39").Value = '7.626'
$ws.Range("E39").Value = '  -2.94%  '

$ws.Range("D40").Value = '0.5534'
$ws.Range("E40").Value = '  -1.84%  '

$ws.Range("D41").Value = '0.1760'
$ws.Range("E41").Value = '  -0.65%  '

$ws.Range("D42").Value = '9.572'
$ws.Range("E42").Value = '  -3.63%  '

$ws.Range("D43").Value = '0.07106'
$ws.Range("E43").Value = '  -3.04%  '

$ws.Range("D44").Value = '11.51'
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").Value = '0.5219'
$ws.Range("E45").Value = '  -1.70%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.129'
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '2.107'
$ws.Range("E47").Value = '  -3.01%  '

$ws.Range("D48").Value = '1.802'
$ws.Range("E48").Value = '  -2.31%  '

$ws.Range("D49").Value = '111.55'
$ws.Range("E49").Value = '  -2.15%  '

$ws.Range("D50").Value = '2.435'
$ws.Range("E50").Value = '  +3.73%  '

$ws.Range("D51").Value = '1.008'
$ws.Range("E51").Value = '  +0.12%  '
